$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log row (row 7) mirroring the existing rows' layout
$row = 7

$ws.Cells.Item(6, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item($row, 1).Value = 42611.887025462966

$ws.Cells.Item($row, 2).Value = 17
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "Random"
